$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final sorted data (descending by value), rows 2-21
$data = @(
    @("English", 26.38367808872994),
    @("Japanese", 8.301906323288847),
    @("Spanish", 8.102184878351892),
    @("German", 6.693468014139586),
    @("Chinese", 5.516175460492082),
    @("Russian", 4.966346720249335),
    @("Arabic", 4.840178086356875),
    @("French", 3.860909789937425),
    @("Italian", 3.741935429607472),
    @("Portuguese", 3.729514418675962),
    @("Malay-Indonesian", 2.434539169680501),
    @("Dutch", 1.68769847709094),
    @("Persian", 1.566519274675442),
    @("Turkish", 1.350213452144083),
    @("Korean", 1.208637958501267),
    @("Thai", 0.9042706229183542),
    @("Urdu", 0.7735871693330135),
    @("Polish", 0.7662396982958973),
    @("Swedish", 0.5346772096167401),
    @("Bengali", 0.3880461993118369)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-obsolete rows 22 and 23 (Uzbek, Vietnamese) entirely
$ws.Range("A22:B23").Delete()
